# Refresh the "cryptos" price/volume snapshot (GitHub Actions scheduled update).
# Price cells in column D are stored as literal text (e.g. "26.167.20", "209.50")
# even though they look numeric, so each is written with a leading apostrophe to
# force text entry and then restyled to "Normal" to avoid leaving the cell on an
# explicit Text number format (matches the original plain/default cell style).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.167.20"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.05%  "
$ws.Range("D3").Value = "'1.579.90"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.33%  "
$ws.Range("E4").Value = "  -0.25%  "
$ws.Range("D5").Value = "'209.50"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.01%  "
$ws.Range("E6").Value = "  -3.17%  "
$ws.Range("E7").Value = "  -0.20%  "
$ws.Range("B8").Value = "Dogecoin"
$ws.Range("C8").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D8").Value = "'0.0609"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.61%  "
$ws.Range("B9").Value = "Cardano"
$ws.Range("C9").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D9").Value = "'0.245"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.69%  "
$ws.Range("D10").Value = "'19.51"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.82%  "
$ws.Range("D11").Value = "'0.0844"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.22%  "
$ws.Range("D12").Value = "'1.802.57"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.33%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "'4.04"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.12%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "'1.569.68"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.95%  "
$ws.Range("E15").Value = "  -1.53%  "
$ws.Range("D16").Value = "'64.44"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.01%  "
$ws.Range("D17").Value = "'26.180.56"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.90%  "
$ws.Range("E18").Value = "  -0.85%  "
$ws.Range("D19").Value = "'7.25"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.21%  "
$ws.Range("E20").Value = "  -0.25%  "
$ws.Range("E21").Value = "  -1.40%  "
$ws.Range("E22").Value = "  -0.55%  "
$ws.Range("E23").Value = "  -2.83%  "
$ws.Range("D24").Value = "'8.88"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.10%  "
$ws.Range("D25").Value = "'144.48"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.56%  "
$ws.Range("E26").Value = "  -0.02%  "
$ws.Range("D27").Value = "'7.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.26%  "
$ws.Range("D28").Value = "'0.112"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.29%  "
$ws.Range("D29").Value = "'15.21"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.15%  "
$ws.Range("E30").Value = "  -1.11%  "
$ws.Range("D31").Value = "'1.14"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.00%  "
$ws.Range("E32").Value = "  -1.85%  "
$ws.Range("E33").Value = "  -0.73%  "
$ws.Range("D34").Value = "'1.276.99"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.14%  "
$ws.Range("D35").Value = "'2.47"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.27%  "
$ws.Range("D36").Value = "'0.610"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.03%  "
$ws.Range("E37").Value = "  -1.13%  "
$ws.Range("E38").Value = "  -2.33%  "
$ws.Range("D39").Value = "'1.09"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -7.18%  "
$ws.Range("D40").Value = "'0.816"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.76%  "
$ws.Range("E41").Value = "  +2.57%  "
$ws.Range("E42").Value = "  -2.84%  "
$ws.Range("E43").Value = "  -2.87%  "
$ws.Range("D44").Value = "'62.42"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.34%  "
$ws.Range("D45").Value = "'1.716.02"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.35%  "
$ws.Range("D46").Value = "'88.98"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.82%  "
$ws.Range("E47").Value = "  -0.12%  "
$ws.Range("E48").Value = "  -2.03%  "
$ws.Range("E49").Value = "  -0.99%  "
$ws.Range("E50").Value = "  -2.17%  "
$ws.Range("E51").Value = "  -0.19%  "
